$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

$ws.Range("C7").Value = "In Progress"
$ws.Range("D7").Value = "☑"
$ws.Range("D7").Select()
